$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "CAOS"
$ws.Range("Q3").Value = "CAOS"
$ws.Range("Q4").Value = "SSK"

$ws.Range("O9").Select()
